$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tableNames = @(
    "domain",
    "vocabulary",
    "redcap_records_tmp_1",
    "care_site",
    "device_exposure",
    "schema_migrations",
    "ar_internal_metadata",
    "redcap_data_dictionaries",
    "concept",
    "concept_relationship",
    "concept_ancestor",
    "drug_strength",
    "drug_era",
    "dose_era",
    "source_to_concept_map",
    "metadata",
    "cdm_source",
    "cohort_definition",
    "person",
    "visit_detail",
    "procedure_occurrence",
    "fact_relationship",
    "location",
    "provider",
    "concept_class",
    "relationship",
    "cohort",
    "cohort_attribute",
    "cost",
    "concept_synonym",
    "attribute_definition",
    "observation_period",
    "death",
    "specimen",
    "visit_occurrence",
    "redcap_variables",
    "redcap_variable_choices",
    "omop_tables",
    "omop_columns",
    "redcap_variable_maps",
    "redcap_variable_choice_maps",
    "redcap_variable_child_maps",
    "redcap_events",
    "redcap_event_maps",
    "redcap_event_map_dependents",
    "redcap_source_links",
    "redcap_projects",
    "drug_exposure",
    "condition_occurrence",
    "measurement",
    "note",
    "payer_plan_period",
    "observation",
    "note_nlp",
    "condition_era"
)

for ($i = 0; $i -lt $tableNames.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $tableNames[$i]
}

# Update the view selection to match the target state as closely as possible
$ws.Range("B1:B55").Select()

Write-Output "Updated $($tableNames.Length) rows in column A"
